# 25 May 2022 Selenium DataDriven Part 1 and Part 2
# Update a few values in the "TestData" sheet and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Replace "Sai" with "Geetha" and mark as Female
$ws.Range("C5").Value = "Geetha"
$ws.Range("H5").Value = "Female"

# Replace "Charan" with "lakshmi" and mark as Female
$ws.Range("C8").Value = "lakshmi"
$ws.Range("H8").Value = "Female"

# Hobbies for row 2 simplified from "Cricket&Hockey" to "Cricket"
$ws.Range("I2").Value = "Cricket"

# Move the active selection to I5
$ws.Activate()
$ws.Range("I5").Select()
